# Updated symbol list on Mon Dec 19 02:59:51 UTC 2022 with GitHub Actions
#
# This script re-applies the "coinranking" symbol-list refresh described by
# the commit: several Price cells get new quotes, a handful of rows
# (8th-17th ranked coins, and the CEJI/BKEXToken pair) shift which coin they
# describe, and a couple of "Bestin24h"/"Worstin24h" suffixes move to a
# different coin's Volume(1h) label.
#
# All of the cells involved in this sheet are stored as text (the Price
# column looks numeric but is kept as a string in the workbook), so plain
# `.Value = "..."` assignment of a numeric-looking string would normally be
# auto-coerced by Excel into a real number. Set-TextValue forces the cell to
# Text format before writing the value and then restores the default
# ("Normal") style so no stray number formats are left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $cellRef, $val) {
    $cell = $sheet.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 - BNB price
Set-TextValue $ws "D2" "250.28"
# Row 3 - OKB price
Set-TextValue $ws "D3" "22.02"
# Row 4 - HuobiToken price
Set-TextValue $ws "D4" "5.525"
# Row 5 - Cronos price
Set-TextValue $ws "D5" "0.05663"
# Row 6 - GateToken price
Set-TextValue $ws "D6" "3.386"
# Row 7 - KuCoinToken price
Set-TextValue $ws "D7" "6.474"
# Row 8 - MXToken price
Set-TextValue $ws "D8" "0.8018"
# Row 9 - FTXToken price
Set-TextValue $ws "D9" "1.048"

# Row 10 - now "One" (was WazirX)
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D10" "0.01170"
$ws.Range("E10").Value = "9OneONEBestin24h"

# Row 11 - now "WazirX" (was MandalaExchangeToken)
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D11" "0.1432"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12 - now "MandalaExchangeToken" (was LiechtensteinCryptoassetsExchange)
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D12" "0.07267"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13 - now "LiechtensteinCryptoassetsExchange" (was BitrueCoin)
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D13" "0.03217"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14 - now "BitrueCoin" (was BitMartToken)
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D14" "0.02945"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15 - now "BitMartToken" (was BitForexToken)
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D15" "0.09264"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16 - now "BitForexToken" (was MCDex)
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D16" "0.001668"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17 - now "MCDex" (was CoinExToken)
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D17" "3.264"
$ws.Range("E17").Value = "16MCDexMCB"

# Row 18 - now "CoinExToken" (was One)
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D18" "0.04744"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19 - TigerCash price
Set-TextValue $ws "D19" "0.006484"

# Row 20 - HotbitToken price, and the "Bestin24h" tag moved off this row
Set-TextValue $ws "D20" "0.005049"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Row 21 - BitKan price
Set-TextValue $ws "D21" "0.001048"

# Row 23 - UpBots price
Set-TextValue $ws "D23" "0.0003206"

# Row 24 - LEO price
Set-TextValue $ws "D24" "4.070"

# Row 25 - BTSEToken price
Set-TextValue $ws "D25" "2.090"

# Row 40 - IDEX price
Set-TextValue $ws "D40" "0.04135"

# Row 41 - KickToken price
Set-TextValue $ws "D41" "0.006913"

# Row 42 - now "CEJI" (was BKEXToken)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.003507"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 - now "BKEXToken" (was CEJI)
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1045"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# Row 44 - LocalTraders price
Set-TextValue $ws "D44" "0.008842"

# Row 45 - CoinLion price
Set-TextValue $ws "D45" "0.00005652"

# Row 47 - CoinbaseStockToken price
Set-TextValue $ws "D47" "0.7867"

# Row 48 - BOLO price
Set-TextValue $ws "D48" "0.01588"
